$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.306.57"
$ws.Range("E2").Value = "  -3.21%  "
$ws.Range("D3").Value = "1.979.60"
$ws.Range("E3").Value = "  -3.92%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.09"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -3.41%  "
$ws.Range("E6").Value = "  -4.55%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "59.10"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -12.55%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  -5.78%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "57.10"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -4.74%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0835"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +8.07%  "
$ws.Range("E12").Value = "  -0.65%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "23.14"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -4.05%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.859"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -8.80%  "
$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "2.271.82"
$ws.Range("E15").Value = "  -3.91%  "
$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.95"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -7.44%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.44"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -4.70%  "
$ws.Range("D18").Value = "1.983.26"
$ws.Range("E18").Value = "  -3.80%  "
$ws.Range("D19").Value = "36.165.03"
$ws.Range("E19").Value = "  -3.50%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "70.36"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -4.69%  "
$ws.Range("D21").Value = "0.0₃0877"
$ws.Range("E21").Value = "  -0.38%  "
$ws.Range("E22").Value = "  -4.23%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "233.71"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -2.97%  "
$ws.Range("E24").Value = "  -0.06%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.52"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -7.13%  "
$ws.Range("E26").Value = "  -6.80%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.89"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -1.98%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "162.73"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.06%  "
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.77"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -2.25%  "
$ws.Range("B30").Value = "Kaspa"
$ws.Range("C30").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.132"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -5.22%  "
$ws.Range("E31").Value = "  -2.61%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.17"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -4.56%  "
$ws.Range("E33").Value = "  -7.38%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0678"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +6.74%  "
$ws.Range("E35").Value = "  -7.48%  "
$ws.Range("E36").Value = "  -1.08%  "
$ws.Range("E37").Value = "  +0.06%  "
$ws.Range("B38").Value = "WEMIXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.82"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.94%  "
$ws.Range("B39").Value = "LidoDAOToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.24"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -8.46%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.96"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -5.42%  "
$ws.Range("E41").Value = "  -4.16%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0956"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -7.71%  "
$ws.Range("E43").Value = "  -5.42%  "
$ws.Range("E44").Value = "  -3.86%  "
$ws.Range("E45").Value = "  -6.04%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "16.13"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -12.54%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "91.97"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -5.93%  "
$ws.Range("D48").Value = "1.363.06"
$ws.Range("E48").Value = "  -4.14%  "
$ws.Range("E49").Value = "  -7.30%  "
$ws.Range("E50").Value = "  -4.27%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "44.92"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -4.90%  "
